# Applies two fleet "operator change" updates:
#  1) GL371TJ: was assigned as "MULETTO FINO AL 31/01/2026 (RICONSEGNARE A MASCIARELLI
#     AL RIENTRO)", now reassigned to CANDELORO AUGUSTA on 2026-01-09.
#  2) GY983FY: was assigned to AUGUSTA.CANDELORO, now returned ("FINE RENT") on 2026-01-09.
# Both changes are reflected in the "Stato Attuale" (current state) sheet and logged
# as updated rows in the "Storico Passaggi" (change history) sheet.
#
# Dates in this workbook are stored as plain text (e.g. "2026-01-09"), not Excel date
# serials. Assigning a date-looking string normally makes Excel auto-convert it to a
# date value, so each date write temporarily forces a Text number format, then restores
# the cell to the default "Normal" style once the literal text value is set.

$wb = $excel.ActiveWorkbook

$wsCurrent = $wb.Worksheets.Item("Stato Attuale")
$wsHistory = $wb.Worksheets.Item("Storico Passaggi")

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# --- Update "Stato Attuale" ---

# Row 51 -> targa GL371TJ
$wsCurrent.Range("B51").Value = "CANDELORO AUGUSTA"
Set-TextValue $wsCurrent.Range("C51") "2026-01-09"

# Row 96 -> targa GY983FY
$wsCurrent.Range("B96").Value = "FINE RENT"
Set-TextValue $wsCurrent.Range("C96") "2026-01-09"

# --- Update "Storico Passaggi" (overwrite the two existing log rows) ---

$wsHistory.Range("A2").Value = "GY983FY"
$wsHistory.Range("B2").Value = "AUGUSTA.CANDELORO"
$wsHistory.Range("C2").Value = "FINE RENT"
Set-TextValue $wsHistory.Range("D2") "2026-01-09"

$wsHistory.Range("A3").Value = "GL371TJ"
$wsHistory.Range("B3").Value = "MULETTO FINO AL 31/01/2026 (RICONSEGNARE A MASCIARELLI AL RIENTRO)"
$wsHistory.Range("C3").Value = "CANDELORO AUGUSTA"
Set-TextValue $wsHistory.Range("D3") "2026-01-09"
